$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---
# These values are plain text in the source data (thousand-separator-dotted
# numbers like "61.199.15", or numeric-looking text like "1.00"/"6.80" whose
# trailing zero must be preserved). Force the cell to text format before
# writing so Excel does not auto-coerce the numeric-looking ones into real
# numbers, then drop back to the default "Normal" style so the cell style
# index matches the original (unstyled) cells.
$dUpdates = @{
    "D2" = "61.199.15"
    "D3" = "2.926.86"
    "D4" = "1.00"
    "D5" = "591.84"
    "D6" = "145.87"
    "D7" = "1.00"
    "D9" = "2.924.29"
    "D10" = "6.80"
    "D11" = "0.145"
    "D12" = "0.444"
    "D13" = "0.0000228"
    "D14" = "33.75"
    "D16" = "3.421.40"
    "D17" = "61.184.29"
    "D18" = "6.75"
    "D19" = "2.928.13"
    "D20" = "431.78"
    "D21" = "13.56"
    "D22" = "0.684"
    "D23" = "7.10"
    "D24" = "80.97"
    "D25" = "10.93"
    "D26" = "2.24"
    "D27" = "12.17"
    "D30" = "1.00"
    "D31" = "2.62"
    "D32" = "7.15"
    "D33" = "26.63"
    "D34" = "0.108"
    "D35" = "0.0₃0867"
    "D36" = "1.01"
    "D37" = "3.12"
    "D38" = "5.64"
    "D39" = "49.83"
    "D40" = "2.02"
    "D41" = "0.124"
    "D42" = "8.63"
    "D43" = "0.291"
    "D44" = "40.12"
    "D45" = "381.55"
    "D46" = "0.0350"
    "D47" = "2.709.74"
    "D48" = "129.99"
    "D50" = "24.25"
}
foreach ($addr in $dUpdates.Keys) {
    $ws.Range($addr).NumberFormat = "@"
}
foreach ($addr in $dUpdates.Keys) {
    $ws.Range($addr).Value = $dUpdates[$addr]
}
foreach ($addr in $dUpdates.Keys) {
    $ws.Range($addr).Style = "Normal"
}

# --- Column E (Volume 1h) updates ---
# Already plain text (percentages padded with two spaces on each side);
# no numeric coercion risk, so these can be written directly.
$eUpdates = @{
    "E2" = "  +0.05%  "
    "E3" = "  -0.99%  "
    "E4" = "  +0.06%  "
    "E5" = "  +0.79%  "
    "E6" = "  -1.51%  "
    "E7" = "  -0.02%  "
    "E8" = "  +1.02%  "
    "E9" = "  -0.49%  "
    "E10" = "  -0.24%  "
    "E11" = "  -0.66%  "
    "E12" = "  -1.37%  "
    "E13" = "  +0.98%  "
    "E14" = "  -2.19%  "
    "E15" = "  +0.31%  "
    "E16" = "  -0.73%  "
    "E17" = "  +0.07%  "
    "E18" = "  -1.92%  "
    "E19" = "  -0.73%  "
    "E20" = "  -0.29%  "
    "E21" = "  -2.42%  "
    "E22" = "  +1.14%  "
    "E23" = "  -2.80%  "
    "E24" = "  +0.61%  "
    "E25" = "  -0.80%  "
    "E26" = "  +1.02%  "
    "E27" = "  +2.20%  "
    "E28" = "  -0.05%  "
    "E29" = "  +7.09%  "
    "E30" = "  +0.17%  "
    "E31" = "  -0.57%  "
    "E32" = "  -3.07%  "
    "E33" = "  -1.16%  "
    "E34" = "  +0.92%  "
    "E35" = "  +3.26%  "
    "E36" = "  -0.05%  "
    "E37" = "  +3.66%  "
    "E38" = "  -1.43%  "
    "E39" = "  -0.56%  "
    "E40" = "  -1.19%  "
    "E41" = "  -1.41%  "
    "E42" = "  -1.54%  "
    "E43" = "  -1.21%  "
    "E44" = "  -5.87%  "
    "E45" = "  +1.80%  "
    "E46" = "  +0.26%  "
    "E47" = "  +1.57%  "
    "E48" = "  -2.64%  "
    "E50" = "  -6.10%  "
    "E51" = "  +0.40%  "
}
foreach ($addr in $eUpdates.Keys) {
    $ws.Range($addr).Value = $eUpdates[$addr]
}
